# Daily attendance processing - 2025-11-27 23:47:30
# Swap the order of names in the "Recorded By" column (G) for rows where
# "dnasr281@gmail.com" is listed first alongside another recorder
# (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()

    if ($val -ne $null -and $val -like "dnasr281@gmail.com, *") {
        $rest = $val.Substring("dnasr281@gmail.com, ".Length)
        $cell.Value = "$rest, dnasr281@gmail.com"
    }
}
